$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptos list refresh (prices + 1h volume deltas), GitHub Actions data pull.
# Column D holds price strings that often look numeric ("1.001", "30.168.01");
# prefix with an apostrophe so Excel keeps them as text like the source data,
# then reset Style to "Normal" so the quote-prefix flag does not linger as a
# visible style change (keeps the cell on the default/unstyled xf, like before).

$ws.Range("D2").Value = "'30.168.01"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.39%  "
$ws.Range("D3").Value = "'1.856.88"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -3.33%  "
$ws.Range("E4").Value = "  +0.22%  "
$ws.Range("D5").Value = "'233.89"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.92%  "
$ws.Range("D6").Value = "'1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.20%  "
$ws.Range("D7").Value = "'0.4674"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.36%  "
$ws.Range("D8").Value = "'0.2815"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.49%  "
$ws.Range("D9").Value = "'0.06546"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.07%  "
$ws.Range("D10").Value = "'20.11"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.77%  "
$ws.Range("D11").Value = "'0.07829"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.59%  "
$ws.Range("D12").Value = "'96.66"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -7.46%  "
$ws.Range("D13").Value = "'1.860.33"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.11%  "
$ws.Range("E14").Value = "  -2.79%  "
$ws.Range("D15").Value = "'0.6648"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.76%  "
$ws.Range("D16").Value = "'282.78"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.97%  "
$ws.Range("D17").Value = "'30.199.48"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.29%  "
$ws.Range("D18").Value = "'1.001"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.16%  "
$ws.Range("D19").Value = "'5.469"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.05%  "
$ws.Range("D20").Value = "'12.59"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.57%  "
$ws.Range("D21").Value = "'2.110.19"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.67%  "
$ws.Range("D22").Value = "'0.000007241"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.82%  "
$ws.Range("D23").Value = "'1.001"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.26%  "
$ws.Range("D24").Value = "'6.138"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.52%  "
$ws.Range("D25").Value = "'167.54"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.04%  "
$ws.Range("D26").Value = "'9.315"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Value = "'18.93"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.69%  "
$ws.Range("D28").Value = "'1.914"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -9.43%  "
$ws.Range("D29").Value = "'1.341"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.61%  "
$ws.Range("D30").Value = "'0.09568"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.39%  "
$ws.Range("D31").Value = "'4.406"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.97%  "
$ws.Range("D32").Value = "'1.469"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.45%  "
$ws.Range("D33").Value = "'4.093"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.66%  "
$ws.Range("D34").Value = "'0.04649"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.32%  "
$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D35").Value = "'1.098"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.85%  "
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").Value = "'0.6992"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.42%  "
$ws.Range("D37").Value = "'1.000"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.30%  "
$ws.Range("E38").Value = "  -0.03%  "
$ws.Range("D39").Value = "'0.01849"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.98%  "
$ws.Range("D40").Value = "'6.397"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.35%  "
$ws.Range("E41").Value = "  -4.45%  "
$ws.Range("D42").Value = "'71.98"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.05%  "
$ws.Range("D43").Value = "'0.8529"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.37%  "
$ws.Range("E44").Value = "  -2.64%  "
$ws.Range("B45").Value = "PaxDollar"
$ws.Range("C45").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D45").Value = "'1.001"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.17%  "
$ws.Range("D46").Value = "'0.4159"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.14%  "
$ws.Range("B47").Value = "Quant"
$ws.Range("C47").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D47").Value = "'103.79"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.41%  "
$ws.Range("D48").Value = "'1.002.07"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.91%  "
$ws.Range("D49").Value = "'7.182"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.05%  "
$ws.Range("D50").Value = "'9.030"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.80%  "
$ws.Range("D51").Value = "'33.89"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.82%  "
